$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 333.1875
$ws.Cells.Item(19, 9).Value = 268.625
$ws.Cells.Item(19, 10).Value = 397.75
$ws.Cells.Item(19, 11).Value = 268.625
$ws.Cells.Item(19, 12).Value = 397.75
$ws.Cells.Item(19, 13).Value = -93.625
$ws.Cells.Item(19, 14).Value = -747.75
# Row 76
$ws.Cells.Item(76, 8).Value = 3093
$ws.Cells.Item(76, 9).Value = 2999.125
$ws.Cells.Item(76, 11).Value = 2999.125
$ws.Cells.Item(76, 13).Value = -2684.125
# Row 79
$ws.Cells.Item(79, 8).Value = 3093
$ws.Cells.Item(79, 9).Value = 2999.125
$ws.Cells.Item(79, 11).Value = 2999.125
$ws.Cells.Item(79, 13).Value = -1907.125
# Row 117
$ws.Cells.Item(117, 8).Value = 48716.8
$ws.Cells.Item(117, 10).Value = 48716.8
$ws.Cells.Item(117, 12).Value = 48716.8
$ws.Cells.Item(117, 14).Value = -57894.8
# Row 130
$ws.Cells.Item(130, 8).Value = 43636
$ws.Cells.Item(130, 10).Value = 43636
$ws.Cells.Item(130, 12).Value = 43636
$ws.Cells.Item(130, 14).Value = -53676
# Row 132
$ws.Cells.Item(132, 8).Value = 14659.214
$ws.Cells.Item(132, 9).Value = 2149.5085
$ws.Cells.Item(132, 10).Value = 81756.73
$ws.Cells.Item(132, 11).Value = 6448.5255
$ws.Cells.Item(132, 12).Value = 245270.19
$ws.Cells.Item(132, 13).Value = -3918.5255
$ws.Cells.Item(132, 14).Value = -250330.19
# Row 137
$ws.Cells.Item(137, 8).Value = 3427.0188
$ws.Cells.Item(137, 9).Value = 1016.5294
$ws.Cells.Item(137, 10).Value = 7740.5264
$ws.Cells.Item(137, 11).Value = 3049.5882
$ws.Cells.Item(137, 12).Value = 23221.5792
$ws.Cells.Item(137, 13).Value = -499.5882000000001
$ws.Cells.Item(137, 14).Value = -28321.5792
# Row 141
$ws.Cells.Item(141, 8).Value = 2105.75
$ws.Cells.Item(141, 9).Value = 870.8823
$ws.Cells.Item(141, 10).Value = 9103.333000000001
$ws.Cells.Item(141, 11).Value = 2612.6469
$ws.Cells.Item(141, 12).Value = 27309.999
$ws.Cells.Item(141, 13).Value = 2567.3531
$ws.Cells.Item(141, 14).Value = -37669.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 1483.591
$ws.Cells.Item(45, 9).Value = 1319.625
$ws.Cells.Item(45, 10).Value = 1920.8334
$ws.Cells.Item(45, 11).Value = 1319.625
$ws.Cells.Item(45, 12).Value = 1920.8334
$ws.Cells.Item(45, 13).Value = -942.625
$ws.Cells.Item(45, 14).Value = -2674.8334
# Row 80
$ws.Cells.Item(80, 8).Value = 37771.89
$ws.Cells.Item(80, 10).Value = 37771.89
$ws.Cells.Item(80, 12).Value = 37771.89
$ws.Cells.Item(80, 14).Value = -39767.89
# Row 83
$ws.Cells.Item(83, 8).Value = 37771.89
$ws.Cells.Item(83, 10).Value = 37771.89
$ws.Cells.Item(83, 12).Value = 113315.67
$ws.Cells.Item(83, 14).Value = -123299.67
# Row 123
$ws.Cells.Item(123, 8).Value = 38496.6
$ws.Cells.Item(123, 10).Value = 38496.6
$ws.Cells.Item(123, 12).Value = 38496.6
$ws.Cells.Item(123, 14).Value = -48296.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 125
$ws.Cells.Item(125, 8).Value = 49768
$ws.Cells.Item(125, 10).Value = 49768
$ws.Cells.Item(125, 12).Value = 49768
$ws.Cells.Item(125, 14).Value = -59608
# Row 126
$ws.Cells.Item(126, 8).Value = 50776
$ws.Cells.Item(126, 10).Value = 50776
$ws.Cells.Item(126, 12).Value = 50776
$ws.Cells.Item(126, 14).Value = -60656
# Row 130
$ws.Cells.Item(130, 8).Value = 49885.5
$ws.Cells.Item(130, 10).Value = 49885.5
$ws.Cells.Item(130, 12).Value = 49885.5
$ws.Cells.Item(130, 14).Value = -59925.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 2220.1667
$ws.Cells.Item(22, 9).Value = 632
$ws.Cells.Item(22, 10).Value = 3808.3333
$ws.Cells.Item(22, 11).Value = 632
$ws.Cells.Item(22, 12).Value = 3808.3333
$ws.Cells.Item(22, 13).Value = -282
$ws.Cells.Item(22, 14).Value = -4508.3333
# Row 31
$ws.Cells.Item(31, 8).Value = 2715.48
$ws.Cells.Item(31, 9).Value = 832.5714
$ws.Cells.Item(31, 10).Value = 3447.7222
$ws.Cells.Item(31, 11).Value = 832.5714
$ws.Cells.Item(31, 12).Value = 3447.7222
$ws.Cells.Item(31, 13).Value = -537.5714
$ws.Cells.Item(31, 14).Value = -4037.7222
# Row 34
$ws.Cells.Item(34, 8).Value = 2715.48
$ws.Cells.Item(34, 9).Value = 832.5714
$ws.Cells.Item(34, 10).Value = 3447.7222
$ws.Cells.Item(34, 11).Value = 832.5714
$ws.Cells.Item(34, 12).Value = 3447.7222
$ws.Cells.Item(34, 13).Value = -630.5714
$ws.Cells.Item(34, 14).Value = -3851.7222
# Row 62
$ws.Cells.Item(62, 8).Value = 3412.7693
$ws.Cells.Item(62, 9).Value = 2928.889
$ws.Cells.Item(62, 11).Value = 2928.889
$ws.Cells.Item(62, 13).Value = -2304.889
# Row 65
$ws.Cells.Item(65, 8).Value = 3412.7693
$ws.Cells.Item(65, 9).Value = 2928.889
$ws.Cells.Item(65, 11).Value = 14644.445
$ws.Cells.Item(65, 13).Value = -11524.445
# Row 100
$ws.Cells.Item(100, 8).Value = 43836
$ws.Cells.Item(100, 10).Value = 43836
$ws.Cells.Item(100, 12).Value = 43836
$ws.Cells.Item(100, 14).Value = -46000
# Row 132
$ws.Cells.Item(132, 8).Value = 45819.625
$ws.Cells.Item(132, 9).Value = 1523.174
$ws.Cells.Item(132, 11).Value = 4569.522
$ws.Cells.Item(132, 13).Value = -2039.522

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 5403.4546
$ws.Cells.Item(5, 9).Value = 6143.9443
$ws.Cells.Item(5, 10).Value = 2071.25
$ws.Cells.Item(5, 11).Value = 18431.8329
$ws.Cells.Item(5, 12).Value = 6213.75
$ws.Cells.Item(5, 13).Value = -18319.8329
$ws.Cells.Item(5, 14).Value = -6437.75
# Row 22
$ws.Cells.Item(22, 8).Value = 9703.833000000001
$ws.Cells.Item(22, 9).Value = 550.5
$ws.Cells.Item(22, 10).Value = 11534.5
$ws.Cells.Item(22, 11).Value = 1651.5
$ws.Cells.Item(22, 12).Value = 34603.5
$ws.Cells.Item(22, 13).Value = -1482.5
$ws.Cells.Item(22, 14).Value = -34941.5
# Row 27
$ws.Cells.Item(27, 8).Value = 9703.833000000001
$ws.Cells.Item(27, 9).Value = 550.5
$ws.Cells.Item(27, 10).Value = 11534.5
$ws.Cells.Item(27, 11).Value = 1651.5
$ws.Cells.Item(27, 12).Value = 34603.5
$ws.Cells.Item(27, 13).Value = -1549.5
$ws.Cells.Item(27, 14).Value = -34807.5
# Row 39
$ws.Cells.Item(39, 8).Value = 3227.2727
$ws.Cells.Item(39, 10).Value = 3227.2727
$ws.Cells.Item(39, 12).Value = 9681.8181
$ws.Cells.Item(39, 14).Value = -10269.8181
# Row 105
$ws.Cells.Item(105, 8).Value = 75785.57000000001
$ws.Cells.Item(105, 9).Value = 2500
$ws.Cells.Item(105, 10).Value = 81422.92
$ws.Cells.Item(105, 11).Value = 7500
$ws.Cells.Item(105, 12).Value = 244268.76
$ws.Cells.Item(105, 13).Value = -4879
$ws.Cells.Item(105, 14).Value = -249510.76
# Row 113
$ws.Cells.Item(113, 8).Value = 2840.913
$ws.Cells.Item(113, 9).Value = 3382.8
$ws.Cells.Item(113, 11).Value = 10148.4
$ws.Cells.Item(113, 13).Value = -7978.400000000001
# Row 135
$ws.Cells.Item(135, 8).Value = 5403.4546
$ws.Cells.Item(135, 9).Value = 6143.9443
$ws.Cells.Item(135, 10).Value = 2071.25
$ws.Cells.Item(135, 11).Value = 55295.4987
$ws.Cells.Item(135, 12).Value = 18641.25
$ws.Cells.Item(135, 13).Value = -52760.4987
$ws.Cells.Item(135, 14).Value = -23711.25
# Row 140
$ws.Cells.Item(140, 8).Value = 192580.81
$ws.Cells.Item(140, 9).Value = 251476.94
$ws.Cells.Item(140, 11).Value = 754430.8200000001
$ws.Cells.Item(140, 13).Value = -749250.8200000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 4768.75
$ws.Cells.Item(80, 9).Value = 5012.5
$ws.Cells.Item(80, 10).Value = 4525
$ws.Cells.Item(80, 11).Value = 5012.5
$ws.Cells.Item(80, 12).Value = 4525
$ws.Cells.Item(80, 13).Value = -4014.5
$ws.Cells.Item(80, 14).Value = -6521
# Row 83
$ws.Cells.Item(83, 8).Value = 4768.75
$ws.Cells.Item(83, 9).Value = 5012.5
$ws.Cells.Item(83, 10).Value = 4525
$ws.Cells.Item(83, 11).Value = 25062.5
$ws.Cells.Item(83, 12).Value = 22625
$ws.Cells.Item(83, 13).Value = -20070.5
$ws.Cells.Item(83, 14).Value = -32609
# Row 102
$ws.Cells.Item(102, 8).Value = 1424.2903
$ws.Cells.Item(102, 9).Value = 1519.7307
$ws.Cells.Item(102, 11).Value = 1519.7307
$ws.Cells.Item(102, 13).Value = 102.2692999999999
# Row 130
$ws.Cells.Item(130, 8).Value = 53986
$ws.Cells.Item(130, 10).Value = 53986
$ws.Cells.Item(130, 12).Value = 53986
$ws.Cells.Item(130, 14).Value = -64026

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 124
$ws.Cells.Item(124, 8).Value = 47496.668
$ws.Cells.Item(124, 10).Value = 47496.668
$ws.Cells.Item(124, 12).Value = 47496.668
$ws.Cells.Item(124, 14).Value = -57316.668
# Row 125
$ws.Cells.Item(125, 8).Value = 48345.5
$ws.Cells.Item(125, 10).Value = 48345.5
$ws.Cells.Item(125, 12).Value = 48345.5
$ws.Cells.Item(125, 14).Value = -58185.5
